$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New gradient rows for the extended 1100-1500m corridor, matching the
# existing distance/MEAN/STD/MIN/MAX/COUNT/Month columns (A:G).
$data = @(
    @(21, 1100, 23.59991455078125, 1.733997464179993, 17.22298049926758, 33.44784545898438, 18476),
    @(22, 1200, 23.59304046630859, 1.765895366668701, 17.20589447021484, 33.2701530456543,  18568),
    @(23, 1300, 23.53413009643555, 1.740623116493225, 17.21273040771484, 31.0728931427002,   18438),
    @(24, 1400, 23.48603439331055, 1.73760187625885,  17.25373649597168, 31.11389923095703,  18366),
    @(25, 1500, 23.4575138092041,  1.805837869644165, 17.40751075744629, 35.96632385253906,  18392)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 1).Value = $entry[1]
    $ws.Cells.Item($row, 2).Value = $entry[2]
    $ws.Cells.Item($row, 3).Value = $entry[3]
    $ws.Cells.Item($row, 4).Value = $entry[4]
    $ws.Cells.Item($row, 5).Value = $entry[5]
    $ws.Cells.Item($row, 6).Value = $entry[6]

    # Column G ("Month") stores the month as text ("10"), matching the
    # existing rows. Assigning the string directly would coerce it to a
    # number, so clone the text cell from the last existing row instead.
    $ws.Range("G20").Copy()
    $ws.Range("G$row").PasteSpecial()
}

$excel.CutCopyMode = $false
